# Add new songs to playlist
#
# - Row 2's song title is corrected from traditional "無名的人" to the
#   simplified-Chinese "无名的人" already used for its author ("毛不易").
# - A brand-new song is inserted right after it: "点燃银河尽头的篝火" by
#   "华晨宇" (YouTube id "yBaYm7Ig7ZQ"), pushing every following row down
#   by one.
# - The old per-row border/centering style (and the now-unused helper
#   column D) is dropped from the playlist rows, matching the rest of
#   the table which uses the sheet's default formatting.
# - The active selection moves to the newly inserted song's YouTube-ID
#   cell (C3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 3+ down to make room for the new song.
$ws.Rows(3).Insert()

# Row 2: same song, but with the simplified-Chinese title.
$ws.Range("A2").Value = "无名的人"
$ws.Range("B2").Value = "毛不易"
$ws.Range("C2").Value = "LgNT-a_ekC8"

# Row 3: the newly added song.
$ws.Range("A3").Value = "点燃银河尽头的篝火"
$ws.Range("B3").Value = "华晨宇"
$ws.Range("C3").Value = "yBaYm7Ig7ZQ"

# These rows (2-5 after the insert) no longer carry the bordered/centered
# style or the spare column D — clean both up.
$ws.Range("A2:C5").ClearFormats()
$ws.Columns("D").Delete()

$ws.Range("C3").Select() | Out-Null
